# Rename the worksheet "Data_Final" to "Data-Final" to avoid loading
# errors when the workbook is opened from Jupyter (e.g. pandas/openpyxl
# sheet-name lookups choking on the underscore).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data_Final")
$ws.Name = "Data-Final"
